{"js": "// no-op\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n"}
